# Adds a "DeleteReason" / "Deleted" column to the single-date report sheets
# (Show, ShowInNewPage, ExportReport, ShowInNewPageDateRange, ExportReportDateRange)
# and widens a couple of columns to fit. Finishes with ExportReportDateRange as
# the active sheet/tab (matching the final selection state of the workbook).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Show" sheet: add column E = DeleteReason / Deleted
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Show")
$ws.Range("E1").Value = "DeleteReason"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E2").Value = "Deleted"
$excel.CutCopyMode = $false
$ws.Range("E1:E2").Select()

# ---------------------------------------------------------------------------
# "ShowInNewPage" sheet: same E column addition
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ShowInNewPage")
$ws.Range("E1").Value = "DeleteReason"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E2").Value = "Deleted"
$excel.CutCopyMode = $false
$ws.Range("E1:E2").Select()

# ---------------------------------------------------------------------------
# "ExportReport" sheet: add E column, and widen D/E a little
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ExportReport")
$ws.Range("E1").Value = "DeleteReason"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E2").Value = "Deleted"
$excel.CutCopyMode = $false
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 11.5
$ws.Range("E1:E2").Select()

# ---------------------------------------------------------------------------
# "ShowInNewPageDateRange" sheet: add column F = DeleteReason / Deleted
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ShowInNewPageDateRange")
$ws.Range("F1").Value = "DeleteReason"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F2").Value = "Deleted"
$excel.CutCopyMode = $false
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Range("E22").Select()

# ---------------------------------------------------------------------------
# "ShowDateRange" sheet: already has the DeleteReason/Deleted columns; just
# widen the new trailing column and move the view/selection.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ShowDateRange")
$ws.Columns.Item(10).ColumnWidth = 11.5
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J1:J2").Select()

# ---------------------------------------------------------------------------
# "ExportReportDateRange" sheet: add column F = DeleteReason / Deleted; this
# ends up the active sheet/tab.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ExportReportDateRange")
$ws.Range("F1").Value = "DeleteReason"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F2").Value = "Deleted"
$excel.CutCopyMode = $false
$ws.Range("F1:F2").Select()
